# Updating filtered feeds from workflow
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: GenomeWeb link for Tempus AI / Verastem story
$ws.Hyperlinks.Add($ws.Range("A8"), "https://www.genomeweb.com/cancer/tempus-ai-verastem-partner-develop-cdx-assay-ovarian-cancer-combination-treatment") | Out-Null
$ws.Range("A8").Style = $ws.Range("A7").Style

# New row 9: 360Dx link for the same Tempus AI / Verastem story
$ws.Hyperlinks.Add($ws.Range("A9"), "https://www.360dx.com/cancer/tempus-ai-verastem-partner-develop-cdx-assay-ovarian-cancer-combination-treatment") | Out-Null
$ws.Range("A9").Style = $ws.Range("A7").Style

$ws.Range("B8").Value = "CDx"
$ws.Range("B9").Value = "CDx"

$ws.Range("C8").Value = "Tempus AI, Verastem Partner to Develop CDx Assay for Ovarian Cancer Combination Treatment"
$ws.Range("C9").Value = "Tempus AI, Verastem Partner to Develop CDx Assay for Ovarian Cancer Combination Treatment"
